$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A44").Value = "ukb51139_subset.csv"
$ws.Range("B44").Value = "28012 x 1081"
$ws.Range("C44").Value = "all"
$ws.Range("D44").Value = "no events"
$ws.Range("E44").Value = "> 140/80"
$ws.Range("F44").Value = "zscore"
$ws.Range("G44").Value = "median"
$ws.Range("H44").Value = "none"
$ws.Range("I44").Value = 25
$ws.Range("K44").Value = "N/A"
$ws.Range("L44").Value = "25.4 & 18.8"
$ws.Range("M44").Value = "33.9 & 29.4"
$ws.Range("N44").Value = "N/A"
$ws.Range("O44").Value = "N/A"
$ws.Range("P44").Value = "filter out large V_sort"

$ws.Range("O44").NumberFormat = $ws.Range("N44").NumberFormat
